# Update Name of Algo
# Apply new KNN imputed values to column B for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12 = 5.492999999999999
    32 = 6.406999999999999
    36 = 8.704000000000001
    38 = 5.438000000000001
    46 = 6.679
    54 = 5.195000000000001
    55 = 4.684
    67 = 5.286
    69 = 5.002
    72 = 5.380999999999999
    91 = 6.378000000000001
    99 = 5.217000000000001
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
